$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Neetcode 150")

# Row 16 (355. Design Twitter) - add notes text first so it lands before the new
# difficulty label in the shared-strings table, then restyle name cell & grow row.
$notes = "1) Twitter class stores a global variable time and users as a dict of User: { userId, tweets: Linked List of { time, tweetId, next }, followers, following } and we store each tweets with latest tweet at the front" + [char]10 + "2) Calculate top 10 tweets by putting the heads of users and following's tweets into a max heap as (-user.tweets.time, user.tweets ) and just pop and reinsert to get the first 10. This way we basically check among all the linked lists simultaneously"

$ws.Range("D16").Value = $notes
$ws.Range("D16").Style = $ws.Range("D15").Style
$ws.Range("D16").WrapText = $true
$ws.Range("D16").VerticalAlignment = $ws.Range("D15").VerticalAlignment

# Row 15 (621. Task Scheduler) - bump difficulty and add "Neutral" highlight style to the name cell
$ws.Range("B15").Value = "Medium (!!!)"
$ws.Range("C15").Style = "Neutral"

# Row 16 (355. Design Twitter) - bump difficulty, restyle name cell, grow row height
$ws.Range("B16").Value = "Medium (!!!)"
$ws.Range("C16").Style = "Neutral"

$ws.Rows.Item(16).RowHeight = 57.6

$ws.Range("D15").Select()
